$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.395.13"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "'2.981.10"
$ws.Range("E3").Value = "  +1.98%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'382.97"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("D6").Value = "'103.00"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("D7").Value = "'0.541"
$ws.Range("E7").Value = "  -1.10%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").Value = "'37.02"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").Value = "'3.450.21"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").Value = "'18.26"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").Value = "'7.58"
$ws.Range("E15").Value = "  +2.61%  "
$ws.Range("D16").Value = "'2.978.20"
$ws.Range("E16").Value = "  +1.91%  "
$ws.Range("E17").Value = "  +7.12%  "
$ws.Range("D18").Value = "'51.370.10"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  -1.47%  "
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").Value = "'12.83"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").Value = "'0.0₃0961"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").Value = "'69.07"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "'261.57"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "'2.91"
$ws.Range("E25").Value = "  +7.74%  "
$ws.Range("D26").Value = "'8.21"
$ws.Range("E26").Value = "  +14.51%  "
$ws.Range("D27").Value = "'7.57"
$ws.Range("E27").Value = "  +10.64%  "
$ws.Range("E28").Value = "  +15.13%  "
$ws.Range("E29").Value = "  -1.05%  "
$ws.Range("D30").Value = "'4.13"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'26.03"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").Value = "'9.88"
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("D34").Value = "'34.72"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'50.98"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").Value = "'2.06"
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("D37").Value = "'0.0454"
$ws.Range("E37").Value = "  +6.58%  "
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("E39").Value = "  -1.24%  "
$ws.Range("D40").Value = "'16.98"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").Value = "'2.59"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("D44").Value = "'122.73"
$ws.Range("E44").Value = "  +2.63%  "
$ws.Range("D45").Value = "'21.72"
$ws.Range("E45").Value = "  -1.41%  "
$ws.Range("D46").Value = "'2.05"
$ws.Range("E46").Value = "  -1.93%  "
$ws.Range("E47").Value = "  +5.81%  "
$ws.Range("E48").Value = "  +2.66%  "
$ws.Range("D49").Value = "'3.30"
$ws.Range("E49").Value = "  +3.17%  "
$ws.Range("D50").Value = "'2.035.21"
$ws.Range("E50").Value = "  +0.13%  "
$ws.Range("E51").Value = "  +3.09%  "
